$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 4731
$ws.Range("F4").Value = 618
$ws.Range("F5").Value = 200
$ws.Range("F6").Value = 1898
$ws.Range("F7").Value = 563
$ws.Range("F8").Value = 787
$ws.Range("F9").Value = 40
$ws.Range("F10").Value = 19
$ws.Range("F11").Value = 421
$ws.Range("F12").Value = 1158
$ws.Range("F13").Value = 1603
$ws.Range("F16").Value = 1907
$ws.Range("F17").Value = 588
$ws.Range("F18").Value = 6
$ws.Range("F21").Value = 217
$ws.Range("F22").Value = 61
$ws.Range("F24").Value = 1204
$ws.Range("F25").Value = 620
$ws.Range("F27").Value = 13
$ws.Range("F29").Value = 298
$ws.Range("F30").Value = 1617
$ws.Range("F32").Value = 503
$ws.Range("F35").Value = 4330

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 35
$ws.Range("F17").Value = 292
$ws.Range("F28").Value = 1747
$ws.Range("F30").Value = 78
$ws.Range("F36").Value = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 100
$ws.Range("F4").Value = 1365
$ws.Range("F5").Value = 1748
$ws.Range("F6").Value = 1094
$ws.Range("F7").Value = 363

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1365
$ws.Range("F4").Value = 1748
$ws.Range("F5").Value = 1094
$ws.Range("F6").Value = 363
$ws.Range("F10").Value = 4731
$ws.Range("F11").Value = 618
$ws.Range("F12").Value = 200
$ws.Range("F13").Value = 1898
$ws.Range("F14").Value = 563
$ws.Range("F15").Value = 787
$ws.Range("F18").Value = 19
$ws.Range("F19").Value = 421
$ws.Range("F20").Value = 1158
$ws.Range("F21").Value = 1603
$ws.Range("F22").Value = 35
$ws.Range("F26").Value = 1907
$ws.Range("F27").Value = 588
$ws.Range("F30").Value = 217
$ws.Range("F32").Value = 61
$ws.Range("F33").Value = 292
$ws.Range("F37").Value = 1204
$ws.Range("F38").Value = 620
$ws.Range("F44").Value = 1617
$ws.Range("F45").Value = 503
$ws.Range("F49").Value = 4330

